$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.08097566666666667
$ws.Range("H2").Value = 0.242927
$ws.Range("I2").Value = 0.005588990034505014
$ws.Range("J2").Value = 0.005588990034505015
$ws.Range("M2").Value = 15.75563966666667
$ws.Range("N2").Value = 47.266919
$ws.Range("O2").Value = 0.3220556913988901
$ws.Range("P2").Value = 0.32205569139889
$ws.Range("Q2").Value = 1.275823425768111
$ws.Range("R2").Value = 11.482410831913
$ws.Range("S2").Value = 0.001799966049784019
$ws.Range("T2").Value = 0.001799966049784019
$ws.Range("G3").Value = 0.08097566666666667
$ws.Range("H3").Value = 0.242927
$ws.Range("I3").Value = 0.005588990034505014
$ws.Range("J3").Value = 0.005588990034505015
$ws.Range("O3").Value = 0.5509544596378365
$ws.Range("P3").Value = 0.5509544596378364
$ws.Range("Q3").Value = 2.182605757048222
$ws.Range("R3").Value = 19.643451813434
$ws.Range("S3").Value = 0.003079278984381963
$ws.Range("T3").Value = 0.003079278984381963
$ws.Range("G4").Value = 0.08097566666666667
$ws.Range("H4").Value = 0.242927
$ws.Range("I4").Value = 0.005588990034505014
$ws.Range("J4").Value = 0.005588990034505015
$ws.Range("O4").Value = 0.1269898489632735
$ws.Range("P4").Value = 0.1269898489632735
$ws.Range("Q4").Value = 0.5030702094981111
$ws.Range("R4").Value = 4.527631885483
$ws.Range("S4").Value = 0.0007097450003390326
$ws.Range("T4").Value = 0.0007097450003390325
$ws.Range("I5").Value = 0.6976944377922635
$ws.Range("J5").Value = 0.6976944377922635
$ws.Range("M5").Value = 15.75563966666667
$ws.Range("N5").Value = 47.266919
$ws.Range("O5").Value = 0.3220556913988901
$ws.Range("P5").Value = 0.32205569139889
$ws.Range("Q5").Value = 159.265789036662
$ws.Range("R5").Value = 1433.392101329958
$ws.Range("S5").Value = 0.2246964645483473
$ws.Range("T5").Value = 0.2246964645483473
$ws.Range("I6").Value = 0.6976944377922635
$ws.Range("J6").Value = 0.6976944377922635
$ws.Range("O6").Value = 0.5509544596378365
$ws.Range("P6").Value = 0.5509544596378364
$ws.Range("S6").Value = 0.3843978619661607
$ws.Range("T6").Value = 0.3843978619661606
$ws.Range("I7").Value = 0.6976944377922635
$ws.Range("J7").Value = 0.6976944377922635
$ws.Range("O7").Value = 0.1269898489632735
$ws.Range("P7").Value = 0.1269898489632735
$ws.Range("S7").Value = 0.08860011127775558
$ws.Range("T7").Value = 0.08860011127775556
$ws.Range("I8").Value = 0.2967165721732315
$ws.Range("J8").Value = 0.2967165721732316
$ws.Range("M8").Value = 15.75563966666667
$ws.Range("N8").Value = 47.266919
$ws.Range("O8").Value = 0.3220556913988901
$ws.Range("P8").Value = 0.32205569139889
$ws.Range("Q8").Value = 67.73280167885467
$ws.Range("R8").Value = 609.595215109692
$ws.Range("S8").Value = 0.09555926080075873
$ws.Range("T8").Value = 0.09555926080075873
$ws.Range("I9").Value = 0.2967165721732315
$ws.Range("J9").Value = 0.2967165721732316
$ws.Range("O9").Value = 0.5509544596378365
$ws.Range("P9").Value = 0.5509544596378364
$ws.Range("S9").Value = 0.1634773186872939
$ws.Range("T9").Value = 0.1634773186872939
$ws.Range("I10").Value = 0.2967165721732315
$ws.Range("J10").Value = 0.2967165721732316
$ws.Range("O10").Value = 0.1269898489632735
$ws.Range("P10").Value = 0.1269898489632735
$ws.Range("S10").Value = 0.03767999268517892
$ws.Range("T10").Value = 0.03767999268517892
